$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 96-97, shifting the existing rows (96-147) down to (98-149).
$ws.Rows("96:97").Insert()

# Fill the newly inserted row 96 with its data.
$ws.Range("A96").Value = 10
$ws.Range("B96").Value = "Vega Modelo de Temuco"
$ws.Range("C96").Value = "La Araucanía"
$ws.Range("D96").Value = "2021-10-19"
$ws.Range("E96").Value = 9
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100101
$ws.Range("H96").Value = "Berries"
$ws.Range("I96").Value = 100112025
$ws.Range("J96").Value = "Frutilla"
$ws.Range("K96").Value = "Sin especificar"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 1500
$ws.Range("N96").Value = 10000
$ws.Range("O96").Value = 10000
$ws.Range("P96").Value = 10000
$ws.Range("Q96").Value = "$/bandeja 7 kilos"
$ws.Range("R96").Value = "Provincia de Melipilla"
$ws.Range("S96").Value = 1429
$ws.Range("T96").Value = 7

# Fill the newly inserted row 97 with its data.
$ws.Range("A97").Value = 10
$ws.Range("B97").Value = "Vega Modelo de Temuco"
$ws.Range("C97").Value = "La Araucanía"
$ws.Range("D97").Value = "2021-10-19"
$ws.Range("E97").Value = 9
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100101
$ws.Range("H97").Value = "Berries"
$ws.Range("I97").Value = 100112025
$ws.Range("J97").Value = "Frutilla"
$ws.Range("K97").Value = "Sin especificar"
$ws.Range("L97").Value = "Tercera"
$ws.Range("M97").Value = 150
$ws.Range("N97").Value = 7000
$ws.Range("O97").Value = 7000
$ws.Range("P97").Value = 7000
$ws.Range("Q97").Value = "$/bandeja 7 kilos"
$ws.Range("R97").Value = "Provincia de Melipilla"
$ws.Range("S97").Value = 1000
$ws.Range("T97").Value = 7
